# Convert the plain-text AKC CHF "pet trusts" URL in the donation table into
# a real hyperlink, matching the existing "AKC Pet Provision Sample" link
# elsewhere in the document (same target URL / relationship, Hyperlink
# character style).

$d = $word.ActiveDocument

$url = "https://www.akcchf.org/how-to-help/donate/planned-giving/pet-trusts.html"

# Locate the plain-text run containing the URL (it is unique in the body).
$rng = $d.Content
$found = $rng.Find.Execute($url, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    # Turn the matched text into a hyperlink pointing at the same address;
    # Word reuses the document's existing relationship for this URL.
    $d.Hyperlinks.Add($rng, $url) | Out-Null
}
